$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to remain text (matches original inline-string cells),
# then clear any style override so no explicit cell style (s attribute) is left behind.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '51.592.42'
$ws.Range("E2").Value = '  -0.15%  '
$ws.Range("D3").Value = '3.107.21'
$ws.Range("E3").Value = '  +2.42%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = '385.75'
$ws.Range("E5").Value = '  +1.66%  '
$ws.Range("D6").Value = '104.03'
$ws.Range("E6").Value = '  +0.66%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").Value = '0.587'
$ws.Range("E9").Value = '  -1.36%  '
$ws.Range("D10").Value = '37.21'
$ws.Range("E10").Value = '  +1.23%  '
$ws.Range("E11").Value = '  +0.14%  '
$ws.Range("D12").Value = '0.0858'
$ws.Range("E12").Value = '  -0.34%  '
$ws.Range("D13").Value = '3.599.55'
$ws.Range("E13").Value = '  +2.37%  '
$ws.Range("D14").Value = '18.59'
$ws.Range("E14").Value = '  +0.12%  '
$ws.Range("E15").Value = '  +1.18%  '
$ws.Range("D16").Value = '3.103.02'
$ws.Range("E16").Value = '  +2.55%  '
$ws.Range("D17").Value = '1.00'
$ws.Range("E17").Value = '  +2.13%  '
$ws.Range("D18").Value = '10.88'
$ws.Range("E18").Value = '  +2.55%  '
$ws.Range("D19").Value = '51.620.94'
$ws.Range("E19").Value = '  -0.11%  '
$ws.Range("D20").Value = '3.26'
$ws.Range("E20").Value = '  +6.93%  '
$ws.Range("D21").Value = '12.47'
$ws.Range("E21").Value = '  -0.42%  '
$ws.Range("D22").Value = '0.0₃0964'
$ws.Range("E22").Value = '  +0.10%  '
$ws.Range("E23").Value = '  +0.10%  '
$ws.Range("D24").Value = '266.98'
$ws.Range("E24").Value = '  -0.75%  '
$ws.Range("E25").Value = '  +0.46%  '
$ws.Range("D26").Value = '8.08'
$ws.Range("E26").Value = '  -1.73%  '
$ws.Range("D27").Value = '27.11'
$ws.Range("E27").Value = '  +3.20%  '
$ws.Range("E28").Value = '  -0.07%  '
$ws.Range("D29").Value = '7.18'
$ws.Range("E29").Value = '  -5.85%  '
$ws.Range("E30").Value = '  -3.48%  '
$ws.Range("E31").Value = '  -1.66%  '
$ws.Range("D32").Value = '10.45'
$ws.Range("E32").Value = '  +1.68%  '
$ws.Range("D33").Value = '0.0483'
$ws.Range("E33").Value = '  +7.43%  '
$ws.Range("E34").Value = '  +3.07%  '
$ws.Range("E35").Value = '  +0.50%  '
$ws.Range("D36").Value = '50.03'
$ws.Range("E36").Value = '  -0.93%  '
$ws.Range("E37").Value = '  -0.14%  '
$ws.Range("E38").Value = '  +0.71%  '
$ws.Range("D39").Value = '0.291'
$ws.Range("E39").Value = '  +1.06%  '
$ws.Range("E40").Value = '  +0.56%  '
$ws.Range("D41").Value = '129.23'
$ws.Range("E41").Value = '  +1.51%  '
$ws.Range("D42").Value = '16.62'
$ws.Range("E42").Value = '  -3.19%  '
$ws.Range("E43").Value = '  -0.01%  '
$ws.Range("E44").Value = '  -3.02%  '
$ws.Range("D45").Value = '3.78'
$ws.Range("E45").Value = '  +0.97%  '
$ws.Range("D46").Value = '22.22'
$ws.Range("E46").Value = '  +1.91%  '
$ws.Range("D47").Value = '2.53'
$ws.Range("E47").Value = '  +5.57%  '
$ws.Range("E48").Value = '  -1.20%  '
$ws.Range("D49").Value = '2.072.52'
$ws.Range("E49").Value = '  +1.85%  '
$ws.Range("E50").Value = '  +20.67%  '
$ws.Range("D51").Value = '0.0320'
$ws.Range("E51").Value = '  -0.13%  '

$ws.Range("D2:E51").Style = "Normal"

Write-Host "Applied 80 cell updates"
